# Apply the "Header Section for Reports" edit to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date range in the header (merged C2:D2)
$ws.Range("C2").Value = "(10/12/2025 - 10/23/2025)"

# Update the service-provider names in column A
$ws.Range("A7").Value = "Nicholas Kavoklis"
$ws.Range("A8").Value = "Roumyadeb Karmakar"

# Update the numeric data table (rows 5-9, columns B-D)
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 7

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 7
